$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.160.84'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '3.506.68'
$ws.Range('E3').Value = '  -5.11%  '
$ws.Range('E4').Value = '  -0.10%  '
$s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.32'
$ws.Range('D5').Style = $s
$ws.Range('E5').Value = '  -0.72%  '
$s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.35'
$ws.Range('D6').Style = $s
$ws.Range('E6').Value = '  -4.26%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '3.498.90'
$ws.Range('E8').Value = '  -5.07%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -5.42%  '
$s = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.79'
$ws.Range('D11').Style = $s
$ws.Range('E11').Value = '  +6.82%  '
$s = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.591'
$ws.Range('D12').Style = $s
$ws.Range('E12').Value = '  -3.31%  '
$s = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.73'
$ws.Range('D13').Style = $s
$ws.Range('E13').Value = '  -5.95%  '
$s = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000276'
$ws.Range('D14').Style = $s
$ws.Range('E14').Value = '  -3.62%  '
$s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '667.86'
$ws.Range('D15').Style = $s
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').Value = '4.063.88'
$ws.Range('E16').Value = '  -5.37%  '
$ws.Range('E17').Value = '  -4.03%  '
$ws.Range('D18').Value = '69.121.62'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('D19').Value = '3.507.64'
$ws.Range('E19').Value = '  -5.20%  '
$ws.Range('E20').Value = '  -1.14%  '
$s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.40'
$ws.Range('D21').Style = $s
$ws.Range('E21').Value = '  -3.63%  '
$s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.16'
$ws.Range('D22').Style = $s
$ws.Range('E22').Value = '  -4.25%  '
$s = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.896'
$ws.Range('D23').Style = $s
$ws.Range('E23').Value = '  -4.81%  '
$s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.16'
$ws.Range('D24').Style = $s
$ws.Range('E24').Value = '  -8.71%  '
$s = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.77'
$ws.Range('D25').Style = $s
$ws.Range('E25').Value = '  -4.47%  '
$s = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.85'
$ws.Range('D26').Style = $s
$ws.Range('E26').Value = '  -4.50%  '
$s = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = $s
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -6.62%  '
$s = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.41'
$ws.Range('D29').Style = $s
$ws.Range('E29').Value = '  -8.40%  '
$s = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.85'
$ws.Range('D30').Style = $s
$ws.Range('E30').Value = '  -7.29%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$s = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.21'
$ws.Range('D31').Style = $s
$ws.Range('E31').Value = '  -6.85%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$s = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.64'
$ws.Range('D32').Style = $s
$ws.Range('E32').Value = '  -6.21%  '
$s = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.34'
$ws.Range('D33').Style = $s
$ws.Range('E33').Value = '  -7.64%  '
$s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.19'
$ws.Range('D34').Style = $s
$ws.Range('E34').Value = '  -2.42%  '
$s = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '613.71'
$ws.Range('D35').Style = $s
$ws.Range('E35').Value = '  +6.75%  '
$s = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.59'
$ws.Range('D36').Style = $s
$ws.Range('E36').Value = '  -13.38%  '
$s = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.86'
$ws.Range('D37').Style = $s
$ws.Range('E37').Value = '  -3.48%  '
$s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.104'
$ws.Range('D38').Style = $s
$ws.Range('E38').Value = '  -4.46%  '
$s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.85'
$ws.Range('D39').Style = $s
$ws.Range('E39').Value = '  -4.77%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$s = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.138'
$ws.Range('D41').Style = $s
$ws.Range('E41').Value = '  -5.89%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0439'
$ws.Range('D42').Style = $s
$ws.Range('E42').Value = '  -5.56%  '
$s = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.334'
$ws.Range('D43').Style = $s
$ws.Range('E43').Value = '  -4.67%  '
$ws.Range('D44').Value = '3.402.89'
$ws.Range('E44').Value = '  -8.67%  '
$s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '33.21'
$ws.Range('D45').Style = $s
$ws.Range('E45').Value = '  -6.72%  '
$ws.Range('D46').Value = '0.0₃0710'
$ws.Range('E46').Value = '  -7.91%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$s = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.59'
$ws.Range('D47').Style = $s
$ws.Range('E47').Value = '  -6.83%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$s = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.83'
$ws.Range('D48').Style = $s
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('E49').Value = '  -1.25%  '
$s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.96'
$ws.Range('D50').Style = $s
$ws.Range('E50').Value = '  -4.52%  '
$s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.66'
$ws.Range('D51').Style = $s
$ws.Range('E51').Value = '  +14.06%  '
